# Apply the StructureDefinition-covered-text.xlsx update:
#  - Metadata sheet: bump Version, Date, Publisher, replace the duplicate
#    "Contact" row with a new "Jurisdiction" row, and delete the extra
#    duplicate "Contact" row entirely (21 rows -> 20 rows).
#  - Elements sheet: update the root Extension row's Short/Definition text
#    to describe the Covered Text extension specifically.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Remove the first of the two duplicate "Contact" rows (row 10). This
# shifts every row below it up by one, turning the former 21-row table
# into a 20-row table.
$meta.Rows.Item(10).Delete()

# Property value updates on the now-shifted rows.
$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet --------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Root "Extension" element row: Short / Definition columns (K / L).
$elements.Range("K2").Value = "Covered Text"
$elements.Range("L2").Value = "Snippet of covered text used as input to the insight asseessment"
